$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "forget pass page check" - Day 9 (Sprint 45) totals were left blank; fill
# in the test-case counts that were forgotten.
$ws.Range("C39").Value = 7250
$ws.Range("C40").Value = 3100
$ws.Range("C41").Value = 3100

# Reflect the scrolled-to / selected cell after filling the values in.
$excel.ActiveWindow.ScrollRow = 22
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C41").Select()
